$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (text content unchanged)
$ws.Range("A1").Value = "Employee Name"
$ws.Range("B1").Value = "Employee PAN"
$ws.Range("C1").Value = "No. of Shares"
$ws.Range("D1").Value = "Date of Allotment(DD-MM-YYYY)"

# Give the newly added rows (4-5) the same date-cell number format as the
# existing D2:D3 cells *before* writing date values into them, so Excel
# reuses the existing style index instead of minting a brand new one.
$ws.Range("D3").Copy()
$ws.Range("D4:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$allotmentDate = (Get-Date -Year 2020 -Month 4 -Day 16 -Hour 0 -Minute 0 -Second 0).Date

# Row 2
$ws.Range("A2").Value = "Harshad Jadhav"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = $allotmentDate

# Row 3
$ws.Range("A3").Value = "vinita mali"
$ws.Range("B3").Value = "12345fdfd"
$ws.Range("C3").Value = 110
$ws.Range("D3").Value = $allotmentDate

# Row 4 (new)
$ws.Range("A4").Value = "vinita mali"
$ws.Range("B4").Value = 3213113113
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = $allotmentDate

# Row 5 (new)
$ws.Range("A5").Value = "Subuser"
$ws.Range("B5").Value = 1333333333
$ws.Range("C5").Value = 40
$ws.Range("D5").Value = $allotmentDate

# Update selection to D2, matching the saved view state
$ws.Range("D2").Select()
